$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New timing values for "Data Set 0 Timings (Pd)" rows 2-6 (columns F,G,H,I,J,L)
# ---------------------------------------------------------------------------
$pdValues = @{
    2 = @{ F = 12.141;              G = 0.342;               H = 8.234999999999999;  I = 1.717; J = 1.722; L = 0.0009300000000000001 }
    3 = @{ F = 13.079;              G = 0.348;               H = 9.087999999999999;  I = 1.729; J = 1.767; L = 0.00088 }
    4 = @{ F = 13.631;              G = 0.376;               H = 9.420999999999999;  I = 1.833; J = 1.853; L = 0.00108 }
    5 = @{ F = 16.754;              G = 0.448;               H = 12.085;             I = 2.076; J = 1.981; L = 0.00102 }
    6 = @{ F = 15.019;              G = 0.483;               H = 10.445;             I = 2.046; J = 1.905; L = 0.0009700000000000001 }
}

# ---------------------------------------------------------------------------
# New timing values for "Data Set 0 Timings (TD)" rows 2-6 (columns F,G,H,I,J,L)
# ---------------------------------------------------------------------------
$tdValues = @{
    2 = @{ F = 13.2;   G = 0.008999999999999999; H = 0.02;  I = 5.597; J = 5.732; L = 0.00117 }
    3 = @{ F = 13.951; G = 0.011;                 H = 0.02;  I = 5.974; J = 5.815; L = 0.00106 }
    4 = @{ F = 13.277; G = 0.008;                 H = 0.016; I = 5.692; J = 5.658; L = 0.00103 }
    5 = @{ F = 13.067; G = 0.008;                 H = 0.016; I = 5.523; J = 5.716; L = 0.0011 }
    6 = @{ F = 13.976; G = 0.01;                  H = 0.02;  I = 6.134; J = 5.898; L = 0.00108 }
}

function Set-RowValues {
    param($ws, $row, $vals)
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("I$row").Value = $vals.I
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("L$row").Value = $vals.L
}

# Sheet 1: "Data Set 0 Timings (Pd)" rows 2-6
$wsPd = $wb.Worksheets.Item("Data Set 0 Timings (Pd)")
foreach ($row in 2..6) {
    Set-RowValues $wsPd $row $pdValues[$row]
}

# Sheet 2: "Data Set 0 Timings (TD)" rows 2-6
$wsTd = $wb.Worksheets.Item("Data Set 0 Timings (TD)")
foreach ($row in 2..6) {
    Set-RowValues $wsTd $row $tdValues[$row]
}

# Sheet 3: "Data Set 0 Timings (combined)" rows 2-6 (Pd values) and rows 7-11 (TD values)
$wsCombined = $wb.Worksheets.Item("Data Set 0 Timings (combined)")
foreach ($row in 2..6) {
    Set-RowValues $wsCombined $row $pdValues[$row]
}
foreach ($row in 2..6) {
    $combinedRow = $row + 5
    Set-RowValues $wsCombined $combinedRow $tdValues[$row]
}
